$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Đơn sale chính")
$ws2 = $wb.Worksheets.Item("Lương")

# Sheet "Đơn sale chính": the totals row no longer carries a discount ratio
$ws1.Range("M3").Value = 0

# Sheet "Lương": remove the per-item "tại HỆ THỐNG" rows (rows 4-10)
$ws2.Rows("4:10").Delete()

# After the shift above, the "Tổng lương tại HỆ THỐNG" row moved from 35 to 28; drop it too
$ws2.Rows("28:28").Delete()

# Update the figures that changed (extra day worked, higher allowance, new totals)
$ws2.Range("B2").Value = 25
$ws2.Range("B3").Value = 875000
$ws2.Range("B28").Value = 705000
$ws2.Range("B31").Value = 705000
